$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws2.Name = "ProductName"
